$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.658.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.66%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.386.95"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.48%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "506.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.36%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.53%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.26%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.552"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.18%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.388.91"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.57%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0979"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.41%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.150"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.35%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.336"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.84%  "

$ws.Range("E13").Value = "  +0.94%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.814.14"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.37%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "56.625.51"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.81%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.79"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.40%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000133"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.72%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.363.83"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.12%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.56%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.80%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "310.04"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.23%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.23"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.40%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.20%  "

$ws.Range("E24").Value = "  +1.23%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.26"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.73%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.43%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.374"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.39%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.151"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.95%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.34"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.85%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "172.81"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.18%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0729"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.32%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.65"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.12%  "

$ws.Range("E33").Value = "  +2.63%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.85"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.24%  "

$ws.Range("E35").Value = "  +0.13%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.997"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.19%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.85"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.22%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.20"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.06%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.85"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.77%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.65"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.37%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.814"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.75%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.43"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.86%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "132.44"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.69%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.40"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.79%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.84"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.05%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.565"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.46%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0910"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.32%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "248.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.40%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0485"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.29%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0210"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.47%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.22"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.12%  "
